$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select column B and delete it, shifting columns C:N left to B:M
$ws.Range("B1:B1048576").Select()
$ws.Range("B:B").Delete() | Out-Null

# Refresh the AutoFilter range/sort state to match the new (narrower) table
$ws.AutoFilterMode = $false
$ws.Range("A1:M1").AutoFilter() | Out-Null
